$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 265.2
$ws.Range("I33").Value = 307.5
$ws.Range("J33").Value = 96
$ws.Range("K33").Value = 307.5
$ws.Range("L33").Value = 96
$ws.Range("M33").Value = -78.5
$ws.Range("N33").Value = -554

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 634.0833
$ws.Range("I38").Value = 260.9
$ws.Range("J38").Value = 2500
$ws.Range("K38").Value = 782.6999999999999
$ws.Range("L38").Value = 7500
$ws.Range("M38").Value = -410.6999999999999
$ws.Range("N38").Value = -8244

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2427.3635
$ws.Range("I51").Value = 10001
$ws.Range("J51").Value = 1670
$ws.Range("K51").Value = 10001
$ws.Range("L51").Value = 1670
$ws.Range("M51").Value = -9517
$ws.Range("N51").Value = -2638

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 112244.445
$ws.Range("I61").Value = 366.66666
$ws.Range("J61").Value = 336000
$ws.Range("K61").Value = 1099.99998
$ws.Range("L61").Value = 1008000
$ws.Range("M61").Value = -927.9999800000001
$ws.Range("N61").Value = -1008344

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2322.9092
$ws.Range("I98").Value = 2322.9092
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2322.9092
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -824.9092000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H110").Value = 45000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 45000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 45000
$ws.Range("N110").Value = -53180

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2322.9092
$ws.Range("I122").Value = 2322.9092
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6968.7276
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4518.7276

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1031.0625
$ws.Range("I129").Value = 430.16666
$ws.Range("J129").Value = 1391.6
$ws.Range("K129").Value = 1290.49998
$ws.Range("L129").Value = 4174.799999999999
$ws.Range("M129").Value = 3709.50002
$ws.Range("N129").Value = -14174.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1175.2258
$ws.Range("I132").Value = 1065.4916
$ws.Range("J132").Value = 3333.3333
$ws.Range("K132").Value = 3196.4748
$ws.Range("L132").Value = 9999.999899999999
$ws.Range("M132").Value = -666.4748
$ws.Range("N132").Value = -15059.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 58950.555
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 58950.555
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 58950.555
$ws.Range("N133").Value = -69070.55499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 75002450
$ws.Range("I135").Value = 33335832
$ws.Range("J135").Value = 200002300
$ws.Range("K135").Value = 300022488
$ws.Range("L135").Value = 1800020700
$ws.Range("M135").Value = -300019953
$ws.Range("N135").Value = -1800025770

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1099101.9
$ws.Range("I138").Value = 1599.15
$ws.Range("J138").Value = 1447515.4
$ws.Range("K138").Value = 4797.450000000001
$ws.Range("L138").Value = 4342546.199999999
$ws.Range("M138").Value = 342.5499999999993
$ws.Range("N138").Value = -4352826.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 13665
$ws.Range("I37").Value = 7997.5
$ws.Range("J37").Value = 25000
$ws.Range("K37").Value = 7997.5
$ws.Range("L37").Value = 25000
$ws.Range("M37").Value = -7724.5
$ws.Range("N37").Value = -25546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 333368320
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 333368320
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 333368320
$ws.Range("N44").Value = -333369296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 119108.61
$ws.Range("I74").Value = 132827.89
$ws.Range("J74").Value = 26503.5
$ws.Range("K74").Value = 132827.89
$ws.Range("L74").Value = 26503.5
$ws.Range("M74").Value = -131953.89
$ws.Range("N74").Value = -28251.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 119108.61
$ws.Range("I77").Value = 132827.89
$ws.Range("J77").Value = 26503.5
$ws.Range("K77").Value = 664139.4500000001
$ws.Range("L77").Value = 132517.5
$ws.Range("M77").Value = -659771.4500000001
$ws.Range("N77").Value = -141253.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2415.8635
$ws.Range("I102").Value = 2207.842
$ws.Range("J102").Value = 3733.3333
$ws.Range("K102").Value = 2207.842
$ws.Range("L102").Value = 3733.3333
$ws.Range("M102").Value = -585.8420000000001
$ws.Range("N102").Value = -6977.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 115.666664
$ws.Range("I22").Value = 83.5
$ws.Range("J22").Value = 180
$ws.Range("K22").Value = 83.5
$ws.Range("L22").Value = 180
$ws.Range("M22").Value = 89.5
$ws.Range("N22").Value = -526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 19888
$ws.Range("I44").Value = 10064
$ws.Range("J44").Value = 24800
$ws.Range("K44").Value = 10064
$ws.Range("L44").Value = 24800
$ws.Range("M44").Value = -9622
$ws.Range("N44").Value = -25684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 30000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 30000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32372

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 30000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 30000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30812

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H90").Value = 30000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 30000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -101856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 30000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 30000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -32808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 607.2
$ws.Range("I105").Value = 607.2
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 607.2
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1139.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 5799.3
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 5799.3
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 17397.9
$ws.Range("N82").Value = -18209.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 5799.3
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 5799.3
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 17397.9
$ws.Range("N85").Value = -20205.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 4800
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 4800
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 14400
$ws.Range("N95").Value = -18518

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 3675.25
$ws.Range("I96").Value = 1101
$ws.Range("J96").Value = 4533.3335
$ws.Range("K96").Value = 3303
$ws.Range("L96").Value = 13600.0005
$ws.Range("M96").Value = -1244
$ws.Range("N96").Value = -17718.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 1303.4286
$ws.Range("I99").Value = 1303.4286
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3910.2858
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1664.2858
$ws.Range("N99").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 3500
$ws.Range("I102").Value = 3000
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 9000
$ws.Range("L102").Value = 18000
$ws.Range("M102").Value = -6566
$ws.Range("N102").Value = -22868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 2342.2222
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 2342.2222
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 7026.6666
$ws.Range("N104").Value = -12268.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2223116.2
$ws.Range("I107").Value = 8547338
$ws.Range("J107").Value = 1092.3243
$ws.Range("K107").Value = 25642014
$ws.Range("L107").Value = 3276.9729
$ws.Range("M107").Value = -25640094
$ws.Range("N107").Value = -7116.9729

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 1000
$ws.Range("I108").Value = 1000
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 3000
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 1800
$ws.Range("I110").Value = 1800
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 5400
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -1310

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 1626.75
$ws.Range("I111").Value = 1626.75
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 4880.25
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -1813.25
$ws.Range("N111").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 799.9429
$ws.Range("I113").Value = 799.92
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 2399.76
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = -229.7599999999998
$ws.Range("N113").Value = -6740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 712.6667
$ws.Range("I114").Value = 630.5714
$ws.Range("J114").Value = 1000
$ws.Range("K114").Value = 1891.7142
$ws.Range("L114").Value = 3000
$ws.Range("M114").Value = 1362.2858
$ws.Range("N114").Value = -9508

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 529
$ws.Range("I117").Value = 529
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1587
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 1855
$ws.Range("N117").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 7663.2
$ws.Range("I120").Value = 8616.666999999999
$ws.Range("J120").Value = 6233
$ws.Range("K120").Value = 25850.001
$ws.Range("L120").Value = 18699
$ws.Range("M120").Value = -21012.001
$ws.Range("N120").Value = -28375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1356940.8
$ws.Range("I139").Value = 2349772.5
$ws.Range("J139").Value = 3079.1365
$ws.Range("K139").Value = 7049317.5
$ws.Range("L139").Value = 9237.4095
$ws.Range("M139").Value = -7044177.5
$ws.Range("N139").Value = -19517.4095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1832000
$ws.Range("I7").Value = 7300000
$ws.Range("J7").Value = 9333.333000000001
$ws.Range("K7").Value = 7300000
$ws.Range("L7").Value = 9333.333000000001
$ws.Range("M7").Value = -7299888
$ws.Range("N7").Value = -9557.333000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 1832000
$ws.Range("I8").Value = 7300000
$ws.Range("J8").Value = 9333.333000000001
$ws.Range("K8").Value = 7300000
$ws.Range("L8").Value = 9333.333000000001
$ws.Range("M8").Value = -7299861
$ws.Range("N8").Value = -9611.333000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 20018
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 20018
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 20018
$ws.Range("M29").Value = ""
$ws.Range("N29").Value = -20608

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2904.3635
$ws.Range("I93").Value = 2294.8
$ws.Range("J93").Value = 9000
$ws.Range("K93").Value = 2294.8
$ws.Range("L93").Value = 9000
$ws.Range("M93").Value = -1046.8
$ws.Range("N93").Value = -11496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 18500
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 18500
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 18500
$ws.Range("N104").Value = -25488

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 79800
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 79800
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 79800
$ws.Range("N112").Value = -82754

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 8166.6665
$ws.Range("I15").Value = 7000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3400.6
$ws.Range("I100").Value = 500
$ws.Range("J100").Value = 5334.3335
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 10668.667
$ws.Range("M100").Value = -459
$ws.Range("N100").Value = -11750.667
